$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: move subject from D2's list encoding into E2, clear D2
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "MCT-3A-Tecnologia da Soldagem"

# Row 3: move subject from D3's list encoding into E3, clear D3
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "MCT-3A-Tecnologia da Soldagem"

# Row 4: clear D4's list encoding
$ws.Range("D4").Value = "-"

# Row 7: move subject from D7's list encoding into F7, clear D7
$ws.Range("D7").Value = "-"
$ws.Range("F7").Value = "MCT-3A-Tecnologia da Soldagem"

# Row 8: set F8 subject
$ws.Range("F8").Value = "MCT-3A-Tecnologia da Soldagem"
